$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 7.833648188735877
$ws.Range("D2").Value = 7.942072044691518
$ws.Range("E2").Value = 13.17348693695376
$ws.Range("F2").Value = 38.80011919111798
$ws.Range("G2").Value = 44.83688766134003
$ws.Range("H2").Value = 18.13860692744838
$ws.Range("J2").Value = 10.14301455620032
$ws.Range("K2").Value = 14.21353379399252
$ws.Range("B3").Value = 7.764242965447843
$ws.Range("D3").Value = 7.904030169677477
$ws.Range("E3").Value = 13.12396650494155
$ws.Range("F3").Value = 38.78703876618852
$ws.Range("G3").Value = 44.76747267378506
$ws.Range("H3").Value = 18.18342773606313
$ws.Range("J3").Value = 10.15163998860108
$ws.Range("K3").Value = 13.85665889925906
$ws.Range("B4").Value = 7.72320144477648
$ws.Range("D4").Value = 7.881715106664447
$ws.Range("E4").Value = 13.09622379112826
$ws.Range("F4").Value = 38.79014149059881
$ws.Range("G4").Value = 44.74055136630349
$ws.Range("H4").Value = 18.21490304301797
$ws.Range("J4").Value = 10.15867910077848
$ws.Range("K4").Value = 13.63553500432987
$ws.Range("B5").Value = 7.706890394296038
$ws.Range("D5").Value = 7.872890157931265
$ws.Range("E5").Value = 13.08559585199511
$ws.Range("F5").Value = 38.79420103139911
$ws.Range("G5").Value = 44.73352497526613
$ws.Range("H5").Value = 18.22872088671235
$ws.Range("J5").Value = 10.16198572518456
$ws.Range("K5").Value = 13.54505852041486
$ws.Range("B6").Value = 7.704207454629304
$ws.Range("D6").Value = 7.871441186785059
$ws.Range("E6").Value = 13.08387221902405
$ws.Range("F6").Value = 38.795043758322
$ws.Range("G6").Value = 44.73259627711145
$ws.Range("H6").Value = 18.23107511202942
$ws.Range("J6").Value = 10.16256124209664
$ws.Range("K6").Value = 13.5300169291296
$ws.Range("B7").Value = 7.722979769947937
$ws.Range("D7").Value = 7.88159499445579
$ws.Range("E7").Value = 13.09607770636398
$ws.Range("F7").Value = 38.79018492934364
$ws.Range("G7").Value = 44.74044064409853
$ws.Range("H7").Value = 18.21508538549696
$ws.Range("J7").Value = 10.15872192142911
$ws.Range("K7").Value = 13.63431610676685
$ws.Range("B8").Value = 7.809401712617622
$ws.Range("D8").Value = 7.928742702994672
$ws.Range("E8").Value = 13.1558644944118
$ws.Range("F8").Value = 38.79329602531135
$ws.Range("G8").Value = 44.80969411935303
$ws.Range("H8").Value = 18.15323839920188
$ws.Range("J8").Value = 10.14562675731575
$ws.Range("K8").Value = 14.0909786682828
$ws.Range("B9").Value = 7.990534219509709
$ws.Range("D9").Value = 8.029179368417871
$ws.Range("E9").Value = 13.29385884648259
$ws.Range("F9").Value = 38.88785393219484
$ws.Range("G9").Value = 45.0700184638651
$ws.Range("H9").Value = 18.06347437382524
$ws.Range("J9").Value = 10.13378207754559
$ws.Range("K9").Value = 14.96494841199359
$ws.Range("B10").Value = 8.129573147518542
$ws.Range("D10").Value = 8.107424063967057
$ws.Range("E10").Value = 13.40734119754875
$ws.Range("F10").Value = 39.01126815018412
$ws.Range("G10").Value = 45.33681298103772
$ws.Range("H10").Value = 18.01691544792497
$ws.Range("J10").Value = 10.13351489008168
$ws.Range("K10").Value = 15.58685502829543
$ws.Range("B11").Value = 8.193867757490745
$ws.Range("D11").Value = 8.143894133026123
$ws.Range("E11").Value = 13.46146170441785
$ws.Range("F11").Value = 39.07908147172608
$ws.Range("G11").Value = 45.47440666320976
$ws.Range("H11").Value = 17.99997932969108
$ws.Range("J11").Value = 10.13522266564675
$ws.Range("K11").Value = 15.8640558084707
$ws.Range("B12").Value = 8.218343311559442
$ws.Range("D12").Value = 8.157822241281364
$ws.Range("E12").Value = 13.48230257026349
$ws.Range("F12").Value = 39.10643126364985
$ws.Range("G12").Value = 45.52882060936633
$ws.Range("H12").Value = 17.99417878219852
$ws.Range("J12").Value = 10.13613195556092
$ws.Range("K12").Value = 15.96810383609796
$ws.Range("B13").Value = 8.213066690055618
$ws.Range("D13").Value = 8.154817473923517
$ws.Range("E13").Value = 13.47779890567554
$ws.Range("F13").Value = 39.10046685930114
$ws.Range("G13").Value = 45.51699926519699
$ws.Range("H13").Value = 17.99540074150738
$ws.Range("J13").Value = 10.13592445365401
$ws.Range("K13").Value = 15.94573775421332
$ws.Range("B14").Value = 8.195878923460468
$ws.Range("D14").Value = 8.145037696466721
$ws.Range("E14").Value = 13.46316942178983
$ws.Range("F14").Value = 39.08129813327359
$ws.Range("G14").Value = 45.4788371811326
$ws.Range("H14").Value = 17.99948982006615
$ws.Range("J14").Value = 10.13529221386508
$ws.Range("K14").Value = 15.8726349278715
$ws.Range("B15").Value = 8.185367020718125
$ws.Range("D15").Value = 8.139062377023599
$ws.Range("E15").Value = 13.45425319973254
$ws.Range("F15").Value = 39.06977398753524
$ws.Range("G15").Value = 45.45576188968541
$ws.Range("H15").Value = 18.0020743743334
$ws.Range("J15").Value = 10.13493912997393
$ws.Range("K15").Value = 15.82773442016955
$ws.Range("B16").Value = 8.125390519160364
$ws.Range("D16").Value = 8.1050575836792
$ws.Range("E16").Value = 13.40385348789128
$ws.Range("F16").Value = 39.00707062206118
$ws.Range("G16").Value = 45.32814556403517
$ws.Range("H16").Value = 18.01810788088758
$ws.Range("J16").Value = 10.13344005719944
$ws.Range("K16").Value = 15.56861561514159
$ws.Range("B17").Value = 8.088849127253599
$ws.Range("D17").Value = 8.084415296414782
$ws.Range("E17").Value = 13.37356573236966
$ws.Range("F17").Value = 38.97158856599494
$ws.Range("G17").Value = 45.25399824027379
$ws.Range("H17").Value = 18.02903272613645
$ws.Range("J17").Value = 10.13298867241294
$ws.Range("K17").Value = 15.40811812028148
$ws.Range("B18").Value = 8.067931160523671
$ws.Range("D18").Value = 8.072625376467959
$ws.Range("E18").Value = 13.35638081184352
$ws.Range("F18").Value = 38.95227943929336
$ws.Range("G18").Value = 45.21287976514533
$ws.Range("H18").Value = 18.03571571312066
$ws.Range("J18").Value = 10.13290126845401
$ws.Range("K18").Value = 15.31527248456644
$ws.Range("B19").Value = 8.06086647638635
$ws.Range("D19").Value = 8.068648013919063
$ws.Range("E19").Value = 13.35060316032587
$ws.Range("F19").Value = 38.94593069671275
$ws.Range("G19").Value = 45.19922101468942
$ws.Range("H19").Value = 18.03804695766287
$ws.Range("J19").Value = 10.1329012658439
$ws.Range("K19").Value = 15.28374844346999
$ws.Range("B20").Value = 8.092728861916063
$ws.Range("D20").Value = 8.086604174881279
$ws.Range("E20").Value = 13.37676559995407
$ws.Range("F20").Value = 38.97525197777923
$ws.Range("G20").Value = 45.26173324739959
$ws.Range("H20").Value = 18.02782841221198
$ws.Range("J20").Value = 10.13301890184694
$ws.Range("K20").Value = 15.42525914313208
$ws.Range("B21").Value = 8.200924070526026
$ws.Range("D21").Value = 8.14790712726491
$ws.Range("E21").Value = 13.46745715103287
$ws.Range("F21").Value = 39.08688319510892
$ws.Range("G21").Value = 45.48998382040917
$ws.Range("H21").Value = 17.99827210832986
$ws.Range("J21").Value = 10.13547079569461
$ws.Range("K21").Value = 15.89413278427503
$ws.Range("B22").Value = 8.272375311571272
$ws.Range("D22").Value = 8.188654143212975
$ws.Range("E22").Value = 13.52874343702659
$ws.Range("F22").Value = 39.16957141453022
$ws.Range("G22").Value = 45.65261039050134
$ws.Range("H22").Value = 17.98252820918743
$ws.Range("J22").Value = 10.13860352817234
$ws.Range("K22").Value = 16.19515430066046
$ws.Range("B23").Value = 8.234179949858531
$ws.Range("D23").Value = 8.166847077008626
$ws.Range("E23").Value = 13.49585375485259
$ws.Range("F23").Value = 39.12455203883543
$ws.Range("G23").Value = 45.56459158046257
$ws.Range("H23").Value = 17.99060332324095
$ws.Range("J23").Value = 10.13679169735476
$ws.Range("K23").Value = 16.03502013106305
$ws.Range("B24").Value = 8.090974551378521
$ws.Range("D24").Value = 8.085614340577752
$ws.Range("E24").Value = 13.37531822898169
$ws.Range("F24").Value = 38.9735923535393
$ws.Range("G24").Value = 45.25823154605798
$ws.Range("H24").Value = 18.02837163005702
$ws.Range("J24").Value = 10.13300469901217
$ws.Range("K24").Value = 15.41751146520545
$ws.Range("B25").Value = 7.94040149458952
$ws.Range("D25").Value = 8.001195656027511
$ws.Range("E25").Value = 13.25436025890614
$ws.Range("F25").Value = 38.8527927178516
$ws.Range("G25").Value = 44.98628311815104
$ws.Range("H25").Value = 18.08436477728321
$ws.Range("J25").Value = 10.13550464887693
$ws.Range("K25").Value = 14.73158294249064
